$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, centered, bordered) from H1 to I1:J1
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Set new header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for new columns I (I0) and J (IF), keyed by row number
$data = @(
    @(2,8,8),
    @(3,7,8),
    @(4,8,8),
    @(5,8,8),
    @(6,7,7),
    @(7,7,7),
    @(8,8,8),
    @(9,8,8),
    @(10,9,9),
    @(11,8,8),
    @(12,8,8),
    @(13,7,7),
    @(14,9,9),
    @(15,9,9),
    @(16,9,9),
    @(17,8,9),
    @(18,9,9),
    @(19,8,8),
    @(20,9,9),
    @(21,9,9),
    @(22,10,10),
    @(23,10,10),
    @(24,8,9),
    @(25,8,9),
    @(26,9,9),
    @(27,7,8),
    @(28,9,9),
    @(29,9,9),
    @(30,9,9),
    @(31,6,7),
    @(32,9,9),
    @(33,7,8),
    @(34,8,8),
    @(35,8,9),
    @(36,9,9),
    @(37,9,9),
    @(38,8,8),
    @(39,7,7),
    @(40,7,8),
    @(41,9,9),
    @(42,9,9),
    @(43,7,8),
    @(44,8,9),
    @(45,9,9),
    @(46,7,8),
    @(47,9,9),
    @(48,7,8),
    @(49,8,8),
    @(50,8,8),
    @(51,5,6),
    @(52,8,8),
    @(53,8,9),
    @(54,7,8),
    @(55,8,8),
    @(56,9,9),
    @(57,8,8),
    @(58,7,8),
    @(59,8,8),
    @(60,7,7),
    @(61,9,9),
    @(62,9,9),
    @(63,6,8),
    @(64,7,8),
    @(65,8,8),
    @(66,8,9),
    @(67,8,8),
    @(68,7,8),
    @(69,7,7),
    @(70,8,8),
    @(71,7,7),
    @(72,6,6),
    @(73,6,7),
    @(74,6,7),
    @(75,9,9),
    @(76,4,4),
    @(77,5,5)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($row, 9).Value = $iVal
    $ws.Cells.Item($row, 10).Value = $jVal
}
